$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 19 ("Afficher en info bulles..." tooltip feature, and the
# "Les fleches..." row) are now Done/Submitted/Approved on 15 Decembre 2014,
# same as the other completed rows (5, 16). Copy the existing "done" cell
# format (green fill) from D16:E16 onto D18:E18 and D19:E19, then stamp the
# date text.
$ws.Range("D16:E16").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Range("D19:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D18").Value = "15 Decembre 2014"
$ws.Range("E18").Value = "15 Decembre 2014"
$ws.Range("D19").Value = "15 Decembre 2014"
$ws.Range("E19").Value = "15 Decembre 2014"

$ws.Range("D18").Select()
